$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 10
$ws.Range("I6").Value = 10
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 30
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = ""
$ws.Range("N6").Value = 82
$ws.Range("H41").Value = 564.8333
$ws.Range("I41").Value = 333
$ws.Range("J41").Value = 796.6667
$ws.Range("K41").Value = 333
$ws.Range("L41").Value = 796.6667
$ws.Range("M41").Value = 107
$ws.Range("N41").Value = -1676.6667
$ws.Range("H70").Value = 4446.625
$ws.Range("J70").Value = 4446.625
$ws.Range("L70").Value = 13339.875
$ws.Range("N70").Value = -13879.875
$ws.Range("H73").Value = 4446.625
$ws.Range("J73").Value = 4446.625
$ws.Range("L73").Value = 13339.875
$ws.Range("N73").Value = -15211.875
$ws.Range("H76").Value = 2627.5715
$ws.Range("I76").Value = 2627.5715
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 2627.5715
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = ""
$ws.Range("N76").Value = -2312.5715
$ws.Range("H79").Value = 2627.5715
$ws.Range("I79").Value = 2627.5715
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 2627.5715
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = ""
$ws.Range("N79").Value = -1535.5715
$ws.Range("H95").Value = 9800
$ws.Range("J95").Value = 9800
$ws.Range("L95").Value = 9800
$ws.Range("N95").Value = -15292
$ws.Range("H101").Value = 544.0909
$ws.Range("I101").Value = 540.1
$ws.Range("J101").Value = 584
$ws.Range("K101").Value = 1620.3
$ws.Range("L101").Value = 1752
$ws.Range("M101").Value = 1.699999999999818
$ws.Range("N101").Value = -4996
$ws.Range("H115").Value = 860
$ws.Range("J115").Value = 3000
$ws.Range("L115").Value = 9000
$ws.Range("N115").Value = -12134
$ws.Range("H138").Value = 1929
$ws.Range("I138").Value = 1148.5
$ws.Range("J138").Value = 1999.9546
$ws.Range("K138").Value = 3445.5
$ws.Range("L138").Value = 5999.8638
$ws.Range("M138").Value = 1694.5
$ws.Range("N138").Value = -16279.8638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21564.4
$ws.Range("I32").Value = 21564.4
$ws.Range("K32").Value = 21564.4
$ws.Range("M32").Value = -21277.4
$ws.Range("H45").Value = 2099.4
$ws.Range("I45").Value = 1999.25
$ws.Range("K45").Value = 1999.25
$ws.Range("M45").Value = -1622.25
$ws.Range("H125").Value = 65000
$ws.Range("J125").Value = 65000
$ws.Range("L125").Value = 65000
$ws.Range("N125").Value = -74840

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3940
$ws.Range("I99").Value = 3940
$ws.Range("K99").Value = 3940
$ws.Range("M99").Value = -2442

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3020
$ws.Range("J31").Value = 3003.25
$ws.Range("L31").Value = 3003.25
$ws.Range("N31").Value = -3593.25
$ws.Range("H34").Value = 3020
$ws.Range("J34").Value = 3003.25
$ws.Range("L34").Value = 3003.25
$ws.Range("N34").Value = -3407.25
$ws.Range("H86").Value = 5915.8335
$ws.Range("I86").Value = 5875
$ws.Range("K86").Value = 5875
$ws.Range("M86").Value = -4752
$ws.Range("H89").Value = 5915.8335
$ws.Range("I89").Value = 5875
$ws.Range("K89").Value = 29375
$ws.Range("M89").Value = -23759
$ws.Range("H105").Value = 791.6111
$ws.Range("I105").Value = 748.5333000000001
$ws.Range("K105").Value = 748.5333000000001
$ws.Range("M105").Value = 998.4666999999999
$ws.Range("H134").Value = 4798
$ws.Range("I134").Value = 2247.5
$ws.Range("K134").Value = 6742.5
$ws.Range("M134").Value = -4207.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 3342.1538
$ws.Range("J22").Value = 3412.4167
$ws.Range("L22").Value = 10237.2501
$ws.Range("N22").Value = -10575.2501
$ws.Range("H27").Value = 3342.1538
$ws.Range("J27").Value = 3412.4167
$ws.Range("L27").Value = 10237.2501
$ws.Range("N27").Value = -10441.2501
$ws.Range("H107").Value = 494.1905
$ws.Range("I107").Value = 271.46155
$ws.Range("J107").Value = 856.125
$ws.Range("K107").Value = 814.38465
$ws.Range("L107").Value = 2568.375
$ws.Range("M107").Value = 1105.61535
$ws.Range("N107").Value = -6408.375
$ws.Range("H108").Value = 372.25
$ws.Range("I108").Value = 372.25
$ws.Range("K108").Value = 1116.75
$ws.Range("M108").Value = 1763.25
$ws.Range("H117").Value = 620
$ws.Range("I117").Value = 620
$ws.Range("K117").Value = 1860
$ws.Range("M117").Value = 1582

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 6057
$ws.Range("J101").Value = 6057
$ws.Range("L101").Value = 6057
$ws.Range("N101").Value = -12547
$ws.Range("H102").Value = 1321.4445
$ws.Range("I102").Value = 815.5
$ws.Range("K102").Value = 815.5
$ws.Range("M102").Value = 806.5
$ws.Range("H126").Value = 3166.6667
$ws.Range("I126").Value = 3166.6667
$ws.Range("K126").Value = 9500.000100000001
$ws.Range("M126").Value = -7030.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1757
$ws.Range("I22").Value = 1433.3334
$ws.Range("J22").Value = 1999.75
$ws.Range("K22").Value = 1433.3334
$ws.Range("L22").Value = 1999.75
$ws.Range("M22").Value = -1138.3334
$ws.Range("N22").Value = -2589.75
$ws.Range("H27").Value = 1757
$ws.Range("I27").Value = 1433.3334
$ws.Range("J27").Value = 1999.75
$ws.Range("K27").Value = 1433.3334
$ws.Range("L27").Value = 1999.75
$ws.Range("M27").Value = -1326.3334
$ws.Range("N27").Value = -2213.75
$ws.Range("H61").Value = 4291.643
$ws.Range("I61").Value = 3340.25
$ws.Range("K61").Value = 3340.25
$ws.Range("M61").Value = -3138.25
$ws.Range("H93").Value = 365.33334
$ws.Range("I93").Value = 365.33334
$ws.Range("K93").Value = 365.33334
$ws.Range("M93").Value = 882.66666
$ws.Range("H113").Value = 4291.643
$ws.Range("I113").Value = 3340.25
$ws.Range("K113").Value = 3340.25
$ws.Range("M113").Value = -1170.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 774.3333
$ws.Range("I100").Value = 649.6667
$ws.Range("J100").Value = 899
$ws.Range("K100").Value = 1299.3334
$ws.Range("L100").Value = 1798
$ws.Range("M100").Value = -758.3334
$ws.Range("N100").Value = -2880
$ws.Range("H113").Value = 1602
$ws.Range("I113").Value = 1469.2222
$ws.Range("J113").Value = 2199.5
$ws.Range("K113").Value = 4407.6666
$ws.Range("L113").Value = 6598.5
$ws.Range("M113").Value = -2237.6666
$ws.Range("N113").Value = -10938.5
$ws.Range("H136").Value = 1790.05
$ws.Range("I136").Value = 1755.6111
$ws.Range("K136").Value = 5266.8333
$ws.Range("M136").Value = -2716.8333
